# V2.0 - heatmaply inclusion + split up server and ui
# Add a new "metadata" worksheet (SampleID / metadata1 / metadata2) after the
# existing ANI and MALDI sheets, populate it, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# --- Update the selection on the ANI sheet (it is no longer the active tab) ---
$wsAni = $wb.Worksheets.Item("ANI")
$wsAni.Range("A2:A13").Select() | Out-Null

# --- Create the new "metadata" worksheet at the end of the workbook ---
$wsMeta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMeta.Name = "metadata"

# Header row
$wsMeta.Range("A1").Value = "SampleID"
$wsMeta.Range("B1").Value = "metadata1"
$wsMeta.Range("C1").Value = "metadata2"

$ids   = @("Sample1", "Sample2", "Sample3", "Sample4", "Sample5", "Sample6", "Sample7", "Sample8", "Sample9", "Sample10", "Sample11", "Sample12")
$meta1 = @("water", "water", "soil", "soil", "soil", "water", "clinical", "clinical", "soil", "water", "clinical", "water")
$meta2 = @("Location 1", "Location 1", "Location 1", "Location 1", "Location 2", "Location 2", "Location 3", "Location 3", "Location 1", "Location 2", "Location 1", "Location 2")

# Fill column by column so the shared-string table is built up in the same
# order as the source workbook (SampleID/metadata1/metadata2 header, then all
# of column B, then all of column C).
for ($i = 0; $i -lt $ids.Length; $i++) {
    $wsMeta.Cells.Item($i + 2, 1).Value = $ids[$i]
}
for ($i = 0; $i -lt $meta1.Length; $i++) {
    $wsMeta.Cells.Item($i + 2, 2).Value = $meta1[$i]
}
for ($i = 0; $i -lt $meta2.Length; $i++) {
    $wsMeta.Cells.Item($i + 2, 3).Value = $meta2[$i]
}

# Column widths (best-fit approximation for SampleID/metadata1 columns)
$wsMeta.Columns.Item(1).ColumnWidth = 12
$wsMeta.Columns.Item(2).ColumnWidth = 9.67

# Select G14 and make this the active sheet/tab, matching the saved view state
$wsMeta.Range("G14").Select() | Out-Null
$wsMeta.Activate()
